# "adding create new client" - adds a new "clients" worksheet at the end of
# the workbook, populated with a header row + one sample client row, plus
# hyperlinks on the website/e-mail columns and text-formatted columns for
# the numeric-looking identifiers (vatRegistrationNo / phone / repPhone).

$wb = $excel.ActiveWorkbook

# --- add the new sheet as the last tab ------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "clients"

# --- header row (row 1) -----------------------------------------------------
$ws.Range("A1").Value = "shortName"
$ws.Range("B1").Value = "ArName"
$ws.Range("C1").Value = "EnName"
$ws.Range("D1").Value = "accountType"
$ws.Range("E1").Value = "relationshipType"
$ws.Range("F1").Value = "corporateType"
$ws.Range("G1").Value = "commercialRegistartionNo"
$ws.Range("H1").Value = "unifiedNo"
$ws.Range("I1").Value = "website"
$ws.Range("J1").NumberFormat = "@"
$ws.Range("J1").Value = "vatRegistrationNo"
$ws.Range("K1").Value = "country"
$ws.Range("L1").Value = "city"
$ws.Range("M1").Value = "district"
$ws.Range("N1").Value = "zip"
$ws.Range("O1").Value = "referenceNo"
$ws.Range("P1").Value = "ArAddress"
$ws.Range("Q1").Value = "EnAddress"
$ws.Range("R1").NumberFormat = "@"
$ws.Range("R1").Value = "phone"
$ws.Range("S1").Value = "mail"
$ws.Range("T1").Value = "repArName"
$ws.Range("U1").Value = "repEnName"
$ws.Range("V1").Value = "ArTitle"
$ws.Range("W1").Value = "EnTitle"
$ws.Range("X1").Value = "repEmail"
$ws.Range("Y1").Value = "status"
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "repPhone"

# --- data row (row 2) -------------------------------------------------------
$ws.Range("A2").Value = "Both"
$ws.Range("B2").Value = "منة الله عماد"
$ws.Range("C2").Value = "Menna Emad"
$ws.Range("D2").Value = "حساب"
$ws.Range("E2").Value = "كلاهما"
$ws.Range("F2").Value = "مؤسسة"
$ws.Range("G2").Value = 1111122222
$ws.Range("H2").Value = 5555566666

# website -> hyperlink
$ws.Range("I2").Value = "https://fai-wp.ahadtest.com"
$ws.Hyperlinks.Add($ws.Range("I2"), "https://fai-wp.ahadtest.com") | Out-Null

# vatRegistrationNo -> kept as text (15-digit number)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "123451234512345"

$ws.Range("K2").Value = "مصر"
$ws.Range("L2").Value = "الإسكندرية"
$ws.Range("M2").Value = "سموحه"
$ws.Range("N2").Value = 123
$ws.Range("O2").Value = 12345
$ws.Range("P2").Value = "سموحه"
$ws.Range("Q2").Value = "smouha"

# phone -> kept as text (leading zero)
$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = "0559505553"

# mail -> hyperlink
$ws.Range("S2").Value = "menna@fai.ws"
$ws.Hyperlinks.Add($ws.Range("S2"), "mailto:menna@fai.ws") | Out-Null

$ws.Range("T2").Value = "ماريو نادى"
$ws.Range("U2").Value = "MarioNady"
$ws.Range("V2").Value = "مدير منتج"
$ws.Range("W2").Value = "product manager"

# repEmail -> hyperlink
$ws.Range("X2").Value = "mario@fai.ws"
$ws.Hyperlinks.Add($ws.Range("X2"), "mailto:mario@fai.ws") | Out-Null

$ws.Range("Y2").Value = "مُعتَمَد"

# repPhone -> kept as text (leading zero)
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "0559505554"

# --- column widths (approximate match to the authored sheet) ---------------
$ws.Columns.Item(1).ColumnWidth = 10.25
$ws.Columns.Item(2).ColumnWidth = 12.6
$ws.Columns.Item(3).ColumnWidth = 11.1
$ws.Columns.Item(4).ColumnWidth = 14.25
$ws.Columns.Item(5).ColumnWidth = 16.25
$ws.Columns.Item(6).ColumnWidth = 12.95
$ws.Columns.Item(7).ColumnWidth = 22.6
$ws.Columns.Item(8).ColumnWidth = 10.1
$ws.Columns.Item(9).ColumnWidth = 24.95
$ws.Columns.Item(10).ColumnWidth = 16.95
$ws.Columns.Item(15).ColumnWidth = 12.45
$ws.Columns.Item(16).ColumnWidth = 9.25
$ws.Columns.Item(17).ColumnWidth = 8.6
$ws.Columns.Item(18).ColumnWidth = 10.75
$ws.Columns.Item(19).ColumnWidth = 14.1
$ws.Columns.Item(20).ColumnWidth = 9.45
$ws.Columns.Item(21).ColumnWidth = 9.95
$ws.Columns.Item(23).ColumnWidth = 14.95
$ws.Columns.Item(24).ColumnWidth = 11.1
$ws.Columns.Item(26).ColumnWidth = 10.25

# --- view state: make "clients" the active/selected sheet+cell -------------
$ws.Range("Z14").Select()
$ws.Activate()
